$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Tabal psicologo" table: drop the "cedula" / "Cadula profecional" row
#    content and replace it with a single "las name" label in column B.
# ---------------------------------------------------------------------------
$ws.Range("A14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("B14").Value = "las name"

# ---------------------------------------------------------------------------
# 2. New "Game class" rows at the bottom of the sheet (37-38).
# ---------------------------------------------------------------------------
$ws.Range("A37").Value = "boolean"
$ws.Range("C37").Value = "active"
# "true" is ambiguous (Excel would coerce it to a Boolean) - force text.
$ws.Range("B37").Value = "'true"
$ws.Range("A38").Value = "int cliente"

# Match the boxed-row look used elsewhere (B18 style) but with only the
# left/right edges kept (no top/bottom) for the new B37:C37 cells.
$boxed = $ws.Range("B22")
foreach ($addr in @("B37", "C37")) {
    $cell = $ws.Range($addr)
    $boxed.Copy($cell)
}
foreach ($addr in @("B37", "C37")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(8).LineStyle = -4142
    $cell.Borders.Item(9).LineStyle = -4142
}
# Re-apply the actual text now that Copy() also duplicated B22's value.
$ws.Range("C37").Value = "active"
$ws.Range("B37").Value = "'true"

# ---------------------------------------------------------------------------
# 3. View bits: scroll so row 29 is at the top and the next empty row (39)
#    is selected, matching where a user would land after typing row 38.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Range("A39").Select()
